$wb = $excel.ActiveWorkbook

# 1. Rename header on "Weekly Quantity" sheet
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header on "Monthly Trend" sheet
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add new "PO Forecast" sheet at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the rest of the workbook's sheets.
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$wsWeekly.Range("A1").Copy()
$headerRange.PasteSpecial(-4122)

# Data rows
$data = @(
    @(44976.99999999999, 185, 77.80402153733189, 291.5121692818328),
    @(45004.99999999999, 168, 51.05578268053675, 279.5770242831935),
    @(45018.99999999999, 159, 44.51922555948484, 272.5851319173335),
    @(45025.99999999999, 155, 42.66069306579083, 265.2200899104586),
    @(45053.99999999999, 137, 30.85440204402834, 242.3962248549254),
    @(45060.99999999999, 133, 22.59442668202372, 241.0986492423849),
    @(45067.99999999999, 129, 16.40832301441396, 243.7054026043563),
    @(45074.99999999999, 124, 16.38830339850651, 234.625511824907),
    @(45144.99999999999, 81, -40.16449264501346, 190.7773591723462),
    @(45151.99999999999, 76, -35.14356937650479, 183.1821805480124),
    @(45158.99999999999, 72, -36.68270089824037, 190.1022965792177),
    @(45165.99999999999, 68, -40.37494741777704, 180.3796364708256),
    @(45172.99999999999, 63, -44.89212940311832, 175.7727546742126),
    @(45179.99999999999, 59, -51.32180805657, 163.976052706169),
    @(45186.99999999999, 55, -53.29611141988718, 167.1347366084073),
    @(45193.99999999999, 50, -53.57625133342427, 159.0801183847998),
    @(45200.99999999999, 46, -58.21759659720453, 158.6551129462819)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Restore original active sheet/selection so the workbook-level view
# (activeTab) is unaffected by adding the new sheet.
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
